# fix: unique command names in XLSX - prefix protocol name to each step
# For every worksheet whose name is a "protocol" sheet (i.e. not one of the
# first five reference sheets), prepend the sheet's own name to the value in
# column A for every data row (row 2 downwards), separated by a single space.

$wb = $excel.ActiveWorkbook

# These sheets are reference/meta sheets and must be left untouched.
$excludedSheets = @("EmilyBellJourney", "NRWaves", "PersonalEmilyBell", "PositiveSpin", "ReEngagement")

foreach ($ws in $wb.Worksheets) {
    $sheetName = $ws.Name

    if ($excludedSheets -contains $sheetName) {
        continue
    }

    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $currentValue = $cell.Value2

        if ($null -eq $currentValue) {
            continue
        }

        $text = [string]$currentValue

        if ($text -eq "") {
            continue
        }

        $prefix = "$sheetName "
        if ($text.StartsWith($prefix)) {
            continue
        }

        $cell.Value = "$sheetName $text"
    }
}
